$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.136.00"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "2.247.75"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.34%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "100.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.573"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("E8").Value = "  +0.03%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.542"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.87%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.74"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.41%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0826"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.49"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").Value = "2.588.20"
$ws.Range("E14").Value = "  +0.08%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.847"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("D16").Value = "2.243.90"
$ws.Range("E16").Value = "  -0.17%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.15"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").Value = "43.990.10"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").Value = "  -7.91%  "
$ws.Range("D20").Value = "0.0₃0972"
$ws.Range("E20").Value = "  -0.52%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.41"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.72%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "65.35"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.08"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.81%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "235.26"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.05"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -6.25%  "
$ws.Range("E26").Value = "  +0.13%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.38"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("E28").Value = "  -0.54%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "37.29"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.13"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.59%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "159.13"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.25%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "20.07"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.24%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0844"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.42%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.68"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.32%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.19"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.94"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +6.23%  "
$ws.Range("E38").Value = "  -2.16%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "16.16"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +12.60%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.65"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.32%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.10"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -6.72%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0314"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.739.73"
$ws.Range("E44").Value = "  -5.89%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.197"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.73%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "81.56"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.40%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "74.34"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "5.14"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.45%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "102.17"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("E50").Value = "  +1.61%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "57.38"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "
